$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.242.32'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '3.249.47'
$ws.Range("E3").Value = '  +6.02%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''576.60'
$ws.Range("E5").Value = '  +3.17%  '
$ws.Range("D6").Value = '''153.29'
$ws.Range("E6").Value = '  +7.32%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '3.239.38'
$ws.Range("E8").Value = '  +6.14%  '
$ws.Range("D9").Value = '''0.514'
$ws.Range("E9").Value = '  +4.24%  '
$ws.Range("D10").Value = '''7.03'
$ws.Range("E10").Value = '  +8.11%  '
$ws.Range("D11").Value = '''0.164'
$ws.Range("E11").Value = '  +4.48%  '
$ws.Range("D12").Value = '''0.489'
$ws.Range("E12").Value = '  +3.63%  '
$ws.Range("D13").Value = '''37.70'
$ws.Range("E13").Value = '  +3.60%  '
$ws.Range("D14").Value = '''0.0000235'
$ws.Range("E14").Value = '  +4.99%  '
$ws.Range("D15").Value = '3.773.78'
$ws.Range("E15").Value = '  +5.98%  '
$ws.Range("D16").Value = '''557.47'
$ws.Range("E16").Value = '  +11.90%  '
$ws.Range("D17").Value = '66.314.01'
$ws.Range("E17").Value = '  +2.82%  '
$ws.Range("D18").Value = '3.255.27'
$ws.Range("E18").Value = '  +5.97%  '
$ws.Range("E19").Value = '  +3.14%  '
$ws.Range("D20").Value = '''7.10'
$ws.Range("E20").Value = '  +5.22%  '
$ws.Range("D21").Value = '''14.42'
$ws.Range("E21").Value = '  +4.15%  '
$ws.Range("D22").Value = '''0.742'
$ws.Range("E22").Value = '  +7.01%  '
$ws.Range("D23").Value = '''7.77'
$ws.Range("E23").Value = '  +7.31%  '
$ws.Range("D24").Value = '''13.58'
$ws.Range("E24").Value = '  +6.00%  '
$ws.Range("D25").Value = '''81.96'
$ws.Range("E25").Value = '  +3.21%  '
$ws.Range("D26").Value = '''0.998'
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = '''9.34'
$ws.Range("E27").Value = '  +17.62%  '
$ws.Range("D28").Value = '''2.96'
$ws.Range("E28").Value = '  +6.19%  '
$ws.Range("D29").Value = '''2.23'
$ws.Range("E29").Value = '  +5.24%  '
$ws.Range("D30").Value = '''27.76'
$ws.Range("E30").Value = '  +5.24%  '
$ws.Range("D31").Value = '''2.75'
$ws.Range("E31").Value = '  +2.93%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = '''1.18'
$ws.Range("E33").Value = '  +5.21%  '
$ws.Range("D34").Value = '''564.04'
$ws.Range("E34").Value = '  +9.86%  '
$ws.Range("D35").Value = '''5.74'
$ws.Range("E35").Value = '  +4.19%  '
$ws.Range("D36").Value = '''6.38'
$ws.Range("E36").Value = '  +5.13%  '
$ws.Range("D37").Value = '''0.0459'
$ws.Range("E37").Value = '  +12.57%  '
$ws.Range("D38").Value = '''55.27'
$ws.Range("E38").Value = '  +2.18%  '
$ws.Range("D39").Value = '''3.10'
$ws.Range("E39").Value = '  +14.50%  '
$ws.Range("D40").Value = '''0.0866'
$ws.Range("E40").Value = '  +7.03%  '
$ws.Range("E41").Value = '  +3.94%  '
$ws.Range("D42").Value = '3.153.58'
$ws.Range("E42").Value = '  +7.15%  '
$ws.Range("D43").Value = '''8.61'
$ws.Range("E43").Value = '  +2.07%  '
$ws.Range("D44").Value = '''0.274'
$ws.Range("E44").Value = '  +10.02%  '
$ws.Range("D45").Value = '''2.28'
$ws.Range("E45").Value = '  +5.36%  '
$ws.Range("D46").Value = '''26.44'
$ws.Range("E46").Value = '  +3.31%  '
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("D48").Value = '0.0₃0557'
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").Value = '''124.56'
$ws.Range("E49").Value = '  +3.05%  '
$ws.Range("E50").Value = '  +2.25%  '
$ws.Range("D51").Value = '''2.24'
$ws.Range("E51").Value = '  +7.20%  '
